# Apply "new version iq and urgency" changes:
#  - Metadata sheet: Status draft -> active, Date updated
#  - Merge the two "Include" sheets into a single "Include #0" sheet that
#    lists the full set of LOINC vital-sign concept codes
#  - Drop the old standalone "Include ValueSet #0" sheet

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsValueSet = $wb.Worksheets.Item("Include ValueSet #0")
$wsInclude = $wb.Worksheets.Item("Include #1")

# ---------------------------------------------------------------------
# 1) Update Metadata sheet values
# ---------------------------------------------------------------------
$wsMeta.Range("B6").Value = "active"
$wsMeta.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# ---------------------------------------------------------------------
# 2) Rebuild the "Include #1" worksheet so it holds all concept rows
#    (it already has the two-column layout/styling we want to keep).
# ---------------------------------------------------------------------

# Move the existing concept rows (rows 2-6: 72514-3, 9269-2, 15074-8,
# blank separator, System URI) down to rows 13-17 so we can make room
# for the 11 new concept codes above them.
$wsInclude.Range("A2:B6").Copy($wsInclude.Range("A13:B17"))

# Stamp the formatting (style only) of the existing data row onto the
# rows that will hold the new concept codes, so no new style entries
# get created.
$wsInclude.Range("A2:B2").Copy()
$wsInclude.Range("A2:B12").PasteSpecial(-4122)

# Fill in the 11 new LOINC concept codes in column A (column B stays
# blank, matching the existing rows). Column B still carries leftover
# values from the format-only paste above (it only copies formatting,
# not content), so explicitly blank it out too.
$newCodes = @("9279-1", "8867-4", "2708-6", "8310-5", "8302-2", "29463-7", "39156-5", "85354-9", "8480-6", "8462-4", "8478-0")
$r = 2
foreach ($code in $newCodes) {
    $wsInclude.Range("A" + $r).Value = $code
    $wsInclude.Range("B" + $r).ClearContents()
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Rename the merged sheet and remove the old standalone
#    "Include ValueSet #0" sheet (rename before delete, since sheet
#    references captured before a delete can go stale afterwards).
# ---------------------------------------------------------------------
$wsInclude.Name = "Include #0"
$wsValueSet.Delete() | Out-Null
